$d = $word.ActiveDocument

# Each of these find-texts is unique within the document (exactly one occurrence),
# so a simple whole-document Find/Replace targets the correct cell in each case.
# The replacements are ordered so that a newly-written value is never
# re-matched by a later Find (handles the 42/3=->78/8= and 78/8=->24/4= overlap).

$d.Content.Find.Execute("33÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "31÷2=", 2) | Out-Null
$d.Content.Find.Execute("94÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "89÷8=", 2) | Out-Null
$d.Content.Find.Execute("61÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "87÷6=", 2) | Out-Null
$d.Content.Find.Execute("78÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "24÷4=", 2) | Out-Null
$d.Content.Find.Execute("42÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "78÷8=", 2) | Out-Null
$d.Content.Find.Execute("65÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "84÷9=", 2) | Out-Null
$d.Content.Find.Execute("53÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "57÷2=", 2) | Out-Null
$d.Content.Find.Execute("28÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "48÷5=", 2) | Out-Null
$d.Content.Find.Execute("90÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "13÷9=", 2) | Out-Null
$d.Content.Find.Execute("46÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "52÷4=", 2) | Out-Null
$d.Content.Find.Execute("48÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "42÷8=", 2) | Out-Null
$d.Content.Find.Execute("83÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "53÷5=", 2) | Out-Null
$d.Content.Find.Execute("66÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "69÷2=", 2) | Out-Null
$d.Content.Find.Execute("93÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "63÷6=", 2) | Out-Null
$d.Content.Find.Execute("61÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "21÷6=", 2) | Out-Null
$d.Content.Find.Execute("39÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "45÷6=", 2) | Out-Null
$d.Content.Find.Execute("68÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "19÷4=", 2) | Out-Null
$d.Content.Find.Execute("33÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "51÷2=", 2) | Out-Null
$d.Content.Find.Execute("81÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "53÷8=", 2) | Out-Null
$d.Content.Find.Execute("36÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "16÷9=", 2) | Out-Null
$d.Content.Find.Execute("82÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "57÷5=", 2) | Out-Null
$d.Content.Find.Execute("11÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "73÷3=", 2) | Out-Null
$d.Content.Find.Execute("88÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "90÷6=", 2) | Out-Null
$d.Content.Find.Execute("71÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "41÷2=", 2) | Out-Null
$d.Content.Find.Execute("58÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "28÷3=", 2) | Out-Null
